$wb = $excel.ActiveWorkbook

# The "About" sheet is the first worksheet in the workbook.
$ws = $wb.Worksheets.Item("About")

# Record a "last updated" style date in C1 as a real Excel date serial
# (44307 == 2021-04-21), formatted with the built-in short-date number
# format (numFmtId 14). Setting the format before the value keeps the
# style table from growing an extra (unused) custom numFmt entry.
$ws.Range("C1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").Value = Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
